$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "B" column values (party/group numbers) that changed in this revision
$ws.Range("B2").Value = 0
$ws.Range("B9").Value = 1
$ws.Range("B21").Value = 2
$ws.Range("B32").Value = 1
$ws.Range("B39").Value = 2
$ws.Range("B44").Value = 2
$ws.Range("B49").Value = 2

# Reset the view: scroll back to the top-left cell and select B1
$ws.Range("B1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
